$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and the OKB/dogwifhat row swap)

$ws.Range("D2").Value = '69.541.13'
$ws.Range("E2").Value = '  -0.50%  '

$ws.Range("D3").Value = '3.778.75'
$ws.Range("E3").Value = '  +0.43%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '''615.58'
$ws.Range("E5").Value = '  -0.88%  '

$ws.Range("D6").Value = '''177.45'
$ws.Range("E6").Value = '  -2.57%  '

$ws.Range("D7").Value = '3.776.18'
$ws.Range("E7").Value = '  +0.49%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  -1.72%  '

$ws.Range("E10").Value = '  -1.51%  '

$ws.Range("D11").Value = '''6.46'
$ws.Range("E11").Value = '  +2.65%  '

$ws.Range("D12").Value = '''0.484'
$ws.Range("E12").Value = '  -1.62%  '

$ws.Range("D13").Value = '''39.78'
$ws.Range("E13").Value = '  -3.93%  '

$ws.Range("E14").Value = '  -2.16%  '

$ws.Range("D15").Value = '4.398.31'
$ws.Range("E15").Value = '  +0.37%  '

$ws.Range("D16").Value = '3.766.50'
$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("D17").Value = '69.611.96'
$ws.Range("E17").Value = '  -0.50%  '

$ws.Range("D18").Value = '''7.56'
$ws.Range("E18").Value = '  -0.66%  '

$ws.Range("E19").Value = '  -3.54%  '

$ws.Range("D20").Value = '''509.76'
$ws.Range("E20").Value = '  +0.12%  '

$ws.Range("D21").Value = '''16.62'
$ws.Range("E21").Value = '  -1.01%  '

$ws.Range("D22").Value = '''9.53'
$ws.Range("E22").Value = '  -0.79%  '

$ws.Range("D23").Value = '''0.734'
$ws.Range("E23").Value = '  +0.63%  '

$ws.Range("D24").Value = '''2.48'
$ws.Range("E24").Value = '  -1.80%  '

$ws.Range("D25").Value = '''86.29'
$ws.Range("E25").Value = '  -1.25%  '

$ws.Range("D26").Value = '''12.90'
$ws.Range("E26").Value = '  -2.05%  '

$ws.Range("E27").Value = '  +4.54%  '

$ws.Range("D28").Value = '''10.57'
$ws.Range("E28").Value = '  -5.31%  '

$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("D30").Value = '''3.01'
$ws.Range("E30").Value = '  +3.16%  '

$ws.Range("D31").Value = '''2.52'
$ws.Range("E31").Value = '  -0.81%  '

$ws.Range("D32").Value = '''8.14'
$ws.Range("E32").Value = '  +2.73%  '

$ws.Range("D33").Value = '''31.17'

$ws.Range("E34").Value = '  -0.30%  '

$ws.Range("D35").Value = '''0.998'
$ws.Range("E35").Value = '  -0.11%  '

$ws.Range("E36").Value = '  -1.41%  '

$ws.Range("D37").Value = '''6.14'
$ws.Range("E37").Value = '  -1.12%  '

$ws.Range("E38").Value = '  +6.67%  '

$ws.Range("D39").Value = '''476.61'
$ws.Range("E39").Value = '  +11.05%  '

$ws.Range("E40").Value = '  +0.59%  '

$ws.Range("D41").Value = '''2.07'
$ws.Range("E41").Value = '  -2.73%  '

$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = '''3.00'
$ws.Range("E42").Value = '  +5.38%  '

$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").Value = '''49.76'
$ws.Range("E43").Value = '  -1.06%  '

$ws.Range("D44").Value = '''44.24'
$ws.Range("E44").Value = '  -3.16%  '

$ws.Range("D45").Value = '''8.59'
$ws.Range("E45").Value = '  -1.98%  '

$ws.Range("D46").Value = '2.946.79'
$ws.Range("E46").Value = '  -2.05%  '

$ws.Range("E47").Value = '  -0.99%  '

$ws.Range("D48").Value = '''27.60'
$ws.Range("E48").Value = '  +0.24%  '

$ws.Range("D49").Value = '''139.60'
$ws.Range("E49").Value = '  +1.49%  '

$ws.Range("E50").Value = '  +0.06%  '

$ws.Range("E51").Value = '  -2.49%  '
